# Update cfb_weather.xlsx with Timestamp 2024-12-07T16:21:31.243141
# This script refreshes the scraped game/weather/odds data on both
# the "FBS" and "Other" sheets to the latest snapshot.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FBS")
$ws2 = $wb.Worksheets.Item("Other")

# --- FBS sheet (rows 2-7) ---
$ws1.Cells.Item(2, 1).Value = 'Clemson @ SMU'
$ws1.Cells.Item(2, 2).Value = 'SAT 12/07'
$ws1.Cells.Item(2, 3).Value = '07:01 PM'
$ws1.Cells.Item(2, 4).Value = 'Low'
$ws1.Cells.Item(2, 5).Value = 'N-S'
$ws1.Cells.Item(2, 6).Value = 'Med'
$ws1.Cells.Item(2, 8).Value = -36.2394104
$ws1.Cells.Item(2, 9).Value = 67.27
$ws1.Cells.Item(2, 10).Value = 62.07
$ws1.Cells.Item(2, 11).Value = 10.3
$ws1.Cells.Item(2, 12).Value = 2000
$ws1.Cells.Item(2, 13).Value = 'N'
$ws1.Cells.Item(2, 14).Value = 'NNW'
$ws1.Cells.Item(2, 15).Value = 44.24
$ws1.Cells.Item(2, 16).Value = 6.4
$ws1.Cells.Item(2, 17).Value = 'N'
$ws1.Cells.Item(2, 18).Value = 1.5
$ws1.Cells.Item(2, 19).Value = -1.5
$ws1.Cells.Item(2, 20).Value = 0
$ws1.Cells.Item(2, 21).Value = -3.9
$ws1.Cells.Item(2, 22).Value = '32.8377223, -96.7827859'
$ws1.Cells.Item(2, 23).Value = 54.5
$ws1.Cells.Item(2, 24).Value = -110
$ws1.Cells.Item(2, 25).Value = 55.5
$ws1.Cells.Item(2, 26).Value = -110
$ws1.Cells.Item(2, 27).Value = -2.5
$ws1.Cells.Item(2, 28).Value = -2.5
$ws1.Cells.Item(2, 31).Value = 0.01834862385321101
$ws1.Cells.Item(2, 32).Value = 0
$ws1.Cells.Item(2, 37).Value = '2024-12-07T16:21:31.243141'
$ws1.Cells.Item(3, 1).Value = 'Marshall @ Louisiana'
$ws1.Cells.Item(3, 2).Value = 'SAT 12/07'
$ws1.Cells.Item(3, 3).Value = '06:30 PM'
$ws1.Cells.Item(3, 4).Value = 'Low'
$ws1.Cells.Item(3, 5).Value = 'NW-SE'
$ws1.Cells.Item(3, 6).Value = 'High'
$ws1.Cells.Item(3, 8).Value = -160.176310297
$ws1.Cells.Item(3, 9).Value = 69.20999999999999
$ws1.Cells.Item(3, 10).Value = 57.18
$ws1.Cells.Item(3, 11).Value = 9.9
$ws1.Cells.Item(3, 12).Value = 1971
$ws1.Cells.Item(3, 13).Value = 'SE'
$ws1.Cells.Item(3, 14).Value = 'NNW'
$ws1.Cells.Item(3, 15).Value = 47.63
$ws1.Cells.Item(3, 16).Value = 3.5
$ws1.Cells.Item(3, 17).Value = 'SE'
$ws1.Cells.Item(3, 18).Value = 0.05
$ws1.Cells.Item(3, 19).Value = 0
$ws1.Cells.Item(3, 20).Value = 0
$ws1.Cells.Item(3, 21).Value = -6.4
$ws1.Cells.Item(3, 22).Value = '30.2158434, -92.0417371'
$ws1.Cells.Item(3, 23).Value = 58.5
$ws1.Cells.Item(3, 24).Value = -110
$ws1.Cells.Item(3, 25).Value = 58.5
$ws1.Cells.Item(3, 26).Value = -110
$ws1.Cells.Item(3, 27).Value = -3.5
$ws1.Cells.Item(3, 28).Value = -5.5
$ws1.Cells.Item(3, 31).Value = 0
$ws1.Cells.Item(3, 32).Value = 2
$ws1.Cells.Item(3, 37).Value = '2024-12-07T16:21:31.243141'
$ws1.Cells.Item(4, 1).Value = 'Ohio @ Miami (OH)'
$ws1.Cells.Item(4, 2).Value = 'SAT 12/07'
$ws1.Cells.Item(4, 3).Value = '12:00 PM'
$ws1.Cells.Item(4, 4).Value = 'High'
$ws1.Cells.Item(4, 5).Value = 'N-S'
$ws1.Cells.Item(4, 6).Value = 'High'
$ws1.Cells.Item(4, 7).Value = 'E'
$ws1.Cells.Item(4, 8).Value = 48.03346249999998
$ws1.Cells.Item(4, 9).Value = 54.38
$ws1.Cells.Item(4, 10).Value = 55.44
$ws1.Cells.Item(4, 11).Value = 10.3
$ws1.Cells.Item(4, 12).Value = 1983
$ws1.Cells.Item(4, 13).Value = 'NNE'
$ws1.Cells.Item(4, 14).Value = 'NNE'
$ws1.Cells.Item(4, 15).Value = 39.02
$ws1.Cells.Item(4, 16).Value = 11.4
$ws1.Cells.Item(4, 17).Value = 'NNE'
$ws1.Cells.Item(4, 18).Value = 0
$ws1.Cells.Item(4, 19).Value = 0
$ws1.Cells.Item(4, 20).Value = 0
$ws1.Cells.Item(4, 21).Value = 1.1
$ws1.Cells.Item(4, 22).Value = '39.5197009, -84.7330255'
$ws1.Cells.Item(4, 23).Value = 44.5
$ws1.Cells.Item(4, 24).Value = -105
$ws1.Cells.Item(4, 25).Value = 44.5
$ws1.Cells.Item(4, 26).Value = -115
$ws1.Cells.Item(4, 31).Value = 0
$ws1.Cells.Item(4, 37).Value = '2024-12-07T16:21:31.243141'
$ws1.Cells.Item(5, 1).Value = 'UNLV @ Boise State'
$ws1.Cells.Item(5, 2).Value = 'FRI 12/06'
$ws1.Cells.Item(5, 3).Value = '06:00 PM'
$ws1.Cells.Item(5, 4).Value = 'High'
$ws1.Cells.Item(5, 5).Value = 'N-S'
$ws1.Cells.Item(5, 6).Value = 'Med'
$ws1.Cells.Item(5, 7).Value = 'E'
$ws1.Cells.Item(5, 9).Value = 53.65
$ws1.Cells.Item(5, 10).Value = 70.04000000000001
$ws1.Cells.Item(5, 11).Value = 6.8
$ws1.Cells.Item(5, 12).Value = 1970
$ws1.Cells.Item(5, 13).Value = 'NNW'
$ws1.Cells.Item(5, 14).Value = 'SSE'
$ws1.Cells.Item(5, 15).Value = 27.02
$ws1.Cells.Item(5, 16).Value = 1.3
$ws1.Cells.Item(5, 17).Value = 'SSE'
$ws1.Cells.Item(5, 18).Value = 0
$ws1.Cells.Item(5, 19).Value = -0.37
$ws1.Cells.Item(5, 20).Value = -0.62
$ws1.Cells.Item(5, 21).Value = -5.5
$ws1.Cells.Item(5, 22).Value = '43.6028839, -116.1958882'
$ws1.Cells.Item(5, 23).Value = 58.5
$ws1.Cells.Item(5, 24).Value = -110
$ws1.Cells.Item(5, 25).Value = 46.5
$ws1.Cells.Item(5, 26).Value = -128
$ws1.Cells.Item(5, 27).Value = -4
$ws1.Cells.Item(5, 28).Value = -4.5
$ws1.Cells.Item(5, 31).Value = -0.2051282051282051
$ws1.Cells.Item(5, 32).Value = 0.5
$ws1.Cells.Item(5, 37).Value = '2024-12-07T16:21:31.243141'
$ws1.Cells.Item(6, 1).Value = 'Western Kentucky @ Jacksonville State'
$ws1.Cells.Item(6, 2).Value = 'FRI 12/06'
$ws1.Cells.Item(6, 3).Value = '06:00 PM'
$ws1.Cells.Item(6, 4).Value = 'Low'
$ws1.Cells.Item(6, 5).Value = 'E-W'
$ws1.Cells.Item(6, 6).Value = 'High'
$ws1.Cells.Item(6, 7).Value = 'N'
$ws1.Cells.Item(6, 8).Value = 43.63323969999999
$ws1.Cells.Item(6, 9).Value = 63.15
$ws1.Cells.Item(6, 10).Value = 59.04
$ws1.Cells.Item(6, 11).Value = 4.8
$ws1.Cells.Item(6, 12).Value = 1947
$ws1.Cells.Item(6, 13).Value = 'SE'
$ws1.Cells.Item(6, 14).Value = 'SSE'
$ws1.Cells.Item(6, 15).Value = 33.08000000000001
$ws1.Cells.Item(6, 16).Value = 2.4
$ws1.Cells.Item(6, 17).Value = 'SSE'
$ws1.Cells.Item(6, 18).Value = 0
$ws1.Cells.Item(6, 19).Value = 0
$ws1.Cells.Item(6, 20).Value = 0
$ws1.Cells.Item(6, 21).Value = -2.4
$ws1.Cells.Item(6, 22).Value = '33.8201052, -85.76647'
$ws1.Cells.Item(6, 23).Value = 58.5
$ws1.Cells.Item(6, 24).Value = -110
$ws1.Cells.Item(6, 25).Value = 61.5
$ws1.Cells.Item(6, 26).Value = -108
$ws1.Cells.Item(6, 27).Value = -3.5
$ws1.Cells.Item(6, 28).Value = -4
$ws1.Cells.Item(6, 31).Value = 0.05128205128205128
$ws1.Cells.Item(6, 32).Value = 0.5
$ws1.Cells.Item(6, 37).Value = '2024-12-07T16:21:31.243141'
$ws1.Cells.Item(7, 1).Value = 'Tulane @ Army'
$ws1.Cells.Item(7, 2).Value = 'FRI 12/06'
$ws1.Cells.Item(7, 3).Value = '08:00 PM'
$ws1.Cells.Item(7, 4).Value = 'High'
$ws1.Cells.Item(7, 5).Value = 'N-S'
$ws1.Cells.Item(7, 6).Value = 'High'
$ws1.Cells.Item(7, 8).Value = 99.21277618000001
$ws1.Cells.Item(7, 9).Value = 52.74
$ws1.Cells.Item(7, 10).Value = 70.11
$ws1.Cells.Item(7, 11).Value = 5.8
$ws1.Cells.Item(7, 12).Value = 1924
$ws1.Cells.Item(7, 13).Value = 'ESE'
$ws1.Cells.Item(7, 14).Value = 'E'
$ws1.Cells.Item(7, 15).Value = 30.62
$ws1.Cells.Item(7, 16).Value = 6.8
$ws1.Cells.Item(7, 17).Value = 'ESE'
$ws1.Cells.Item(7, 18).Value = 0
$ws1.Cells.Item(7, 19).Value = 0
$ws1.Cells.Item(7, 20).Value = -0.17
$ws1.Cells.Item(7, 21).Value = 1
$ws1.Cells.Item(7, 22).Value = '41.3874924, -73.9640891'
$ws1.Cells.Item(7, 23).Value = 48.5
$ws1.Cells.Item(7, 24).Value = -114
$ws1.Cells.Item(7, 25).Value = 51.5
$ws1.Cells.Item(7, 26).Value = -118
$ws1.Cells.Item(7, 27).Value = 4
$ws1.Cells.Item(7, 28).Value = 4.5
$ws1.Cells.Item(7, 31).Value = 0.06185567010309279
$ws1.Cells.Item(7, 32).Value = -0.5
$ws1.Cells.Item(7, 37).Value = '2024-12-07T16:21:31.243141'
$ws1.Cells.Item(2, 7).Value = ""
$ws1.Cells.Item(4, 27).Value = ""
$ws1.Cells.Item(4, 28).Value = ""
$ws1.Cells.Item(4, 32).Value = ""
$ws1.Cells.Item(5, 8).Value = ""

# --- Other sheet (rows 2-5) ---
$ws2.Cells.Item(2, 1).Value = 'Montana vs South Dakota State'
$ws2.Cells.Item(2, 2).Value = 'South Dakota State'
$ws2.Cells.Item(2, 3).Value = 'Montana'
$ws2.Cells.Item(2, 4).Value = 'SAT 12/07'
$ws2.Cells.Item(2, 5).Value = '01:00 PM'
$ws2.Cells.Item(2, 6).Value = 'Mid'
$ws2.Cells.Item(2, 10).Value = -474.5684815
$ws2.Cells.Item(2, 11).Value = 46.7
$ws2.Cells.Item(2, 12).Value = 47.64
$ws2.Cells.Item(2, 14).Value = 2016
$ws2.Cells.Item(2, 15).Value = 'E'
$ws2.Cells.Item(2, 16).Value = 'E'
$ws2.Cells.Item(2, 17).Value = 51.38
$ws2.Cells.Item(2, 18).Value = 13.6
$ws2.Cells.Item(2, 19).Value = 'E'
$ws2.Cells.Item(2, 20).Value = 0
$ws2.Cells.Item(2, 21).Value = -2
$ws2.Cells.Item(2, 22).Value = 0
$ws2.Cells.Item(2, 24).Value = '44.3210182, -96.7801386'
$ws2.Cells.Item(3, 1).Value = 'Rhode Island vs Mercer'
$ws2.Cells.Item(3, 2).Value = 'Mercer'
$ws2.Cells.Item(3, 3).Value = 'Rhode Island'
$ws2.Cells.Item(3, 4).Value = 'SAT 12/07'
$ws2.Cells.Item(3, 5).Value = '02:00 PM'
$ws2.Cells.Item(3, 6).Value = 'Low'
$ws2.Cells.Item(3, 10).Value = 105.98195272
$ws2.Cells.Item(3, 11).Value = 64.83
$ws2.Cells.Item(3, 12).Value = 52.81
$ws2.Cells.Item(3, 14).Value = 2013
$ws2.Cells.Item(3, 15).Value = 'NNE'
$ws2.Cells.Item(3, 16).Value = 'NE'
$ws2.Cells.Item(3, 17).Value = 55.7
$ws2.Cells.Item(3, 18).Value = 4
$ws2.Cells.Item(3, 19).Value = 'NE'
$ws2.Cells.Item(3, 20).Value = 0
$ws2.Cells.Item(3, 21).Value = 0
$ws2.Cells.Item(3, 22).Value = 0
$ws2.Cells.Item(3, 24).Value = '32.8262075, -83.6522485'
$ws2.Cells.Item(4, 1).Value = 'Villanova vs Incarnate Word'
$ws2.Cells.Item(4, 2).Value = 'Incarnate Word'
$ws2.Cells.Item(4, 3).Value = 'Villanova'
$ws2.Cells.Item(4, 4).Value = 'SAT 12/07'
$ws2.Cells.Item(4, 5).Value = '01:00 PM'
$ws2.Cells.Item(4, 6).Value = 'Low'
$ws2.Cells.Item(4, 10).Value = 81.08228299999999
$ws2.Cells.Item(4, 11).Value = 70.73999999999999
$ws2.Cells.Item(4, 12).Value = 55.05
$ws2.Cells.Item(4, 14).Value = 2008
$ws2.Cells.Item(4, 15).Value = 'SSW'
$ws2.Cells.Item(4, 16).Value = 'SSW'
$ws2.Cells.Item(4, 17).Value = 45.14000000000001
$ws2.Cells.Item(4, 18).Value = 6.4
$ws2.Cells.Item(4, 19).Value = 'SSW'
$ws2.Cells.Item(4, 20).Value = 0.2
$ws2.Cells.Item(4, 21).Value = 0
$ws2.Cells.Item(4, 22).Value = 0
$ws2.Cells.Item(4, 24).Value = '29.4674787, -98.470014'
$ws2.Cells.Item(5, 1).Value = 'Illinois State vs UC Davis'
$ws2.Cells.Item(5, 2).Value = 'UC Davis'
$ws2.Cells.Item(5, 3).Value = 'Illinois State'
$ws2.Cells.Item(5, 4).Value = 'SAT 12/07'
$ws2.Cells.Item(5, 5).Value = '01:00 PM'
$ws2.Cells.Item(5, 6).Value = 'High'
$ws2.Cells.Item(5, 10).Value = -231.4896765
$ws2.Cells.Item(5, 11).Value = 62.21
$ws2.Cells.Item(5, 12).Value = 53.1
$ws2.Cells.Item(5, 14).Value = 2007
$ws2.Cells.Item(5, 15).Value = 'WSW'
$ws2.Cells.Item(5, 16).Value = 'E'
$ws2.Cells.Item(5, 17).Value = 60.68000000000001
$ws2.Cells.Item(5, 18).Value = 1.4
$ws2.Cells.Item(5, 19).Value = 'E'
$ws2.Cells.Item(5, 20).Value = 0
$ws2.Cells.Item(5, 21).Value = 0
$ws2.Cells.Item(5, 22).Value = 0
$ws2.Cells.Item(5, 24).Value = '38.5365266, -121.7627936'